# Issue 260 update: append a new benchmark column ("v1433") to the two
# performance sheets (Sponza -> sheet2, ComplexMesh -> sheet3).
#
# Sponza gets a new column M (after existing L = v1423).
# ComplexMesh gets a new column L (after existing K = v1423).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Sponza" (sheet2.xml) — add column M
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Sponza")

# Clone number formats / styles from column L into the new column M so the
# new cells pick up the same look (header style, value style, avg style...).
$ws.Range("L1:L16").Copy()
$ws.Range("M1:M16").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Header
$ws.Range("M1").Value = "v1433"

# Raw sample values (rows 2-11)
$ws.Range("M2").Value = 7019
$ws.Range("M3").Value = 7019
$ws.Range("M4").Value = 6970
$ws.Range("M5").Value = 6983
$ws.Range("M6").Value = 6982
$ws.Range("M7").Value = 6956
$ws.Range("M8").Value = 6984
$ws.Range("M9").Value = 6982
$ws.Range("M10").Value = 7001
$ws.Range("M11").Value = 6972

# Row 12: AVG
$ws.Range("M12").Formula = "=AVERAGE(M2:M11)"
# Row 13: VAR
$ws.Range("M13").Formula = "=_xlfn.VAR.S(M2:M11)"
# Row 14: DIFF ACCEPT (vs previous column)
$ws.Range("M14").Formula = "=1-_xlfn.T.TEST(L2:L11,M2:M11,2,3)"
# Row 15: Perf (Step) (vs previous column)
$ws.Range("M15").Formula = "=L12/M12"
# Row 16: Perf (Total) (vs first column)
$ws.Range("M16").Formula = "=B12/M12"

# Dimension / selection
[void]$ws.Range("M2").Select()

# Extend the conditional formatting that highlighted B15:L16 to cover the
# new column.
$fcs = $ws.Range("B15:L16").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("B15:M16"))
}

# ---------------------------------------------------------------------
# Sheet "ComplexMesh" (sheet3.xml) — add column L
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ComplexMesh")

$ws2.Range("K1:K16").Copy()
$ws2.Range("L1:L16").PasteSpecial(-4122)
$ws2.Application.CutCopyMode = $false

# Header
$ws2.Range("L1").Value = "v1433"

# Raw sample values (rows 2-11)
$ws2.Range("L2").Value = 4917
$ws2.Range("L3").Value = 4967
$ws2.Range("L4").Value = 5046
$ws2.Range("L5").Value = 4977
$ws2.Range("L6").Value = 4964
$ws2.Range("L7").Value = 4943
$ws2.Range("L8").Value = 4952
$ws2.Range("L9").Value = 4947
$ws2.Range("L10").Value = 4949
$ws2.Range("L11").Value = 4944

# Row 12: AVG column here is a literal value (matches H12:K12 pattern on
# this sheet, which are plain numbers rather than AVERAGE() formulas).
$ws2.Range("L12").Value = 4964

# Row 13: VAR
$ws2.Range("L13").Formula = "=_xlfn.VAR.S(L2:L11)"
# Row 14: DIFF ACCEPT (vs previous column)
$ws2.Range("L14").Formula = "=1-_xlfn.T.TEST(K2:K11,L2:L11,2,3)"
# Row 15: Perf (Step) (vs previous column)
$ws2.Range("L15").Formula = "=K12/L12"
# Row 16: Perf (Total) (vs first column)
$ws2.Range("L16").Formula = "=B12/L12"

# Dimension / selection
[void]$ws2.Range("L2").Select()

# Extend the conditional formatting that highlighted B15:K16 to cover the
# new column.
$fcs2 = $ws2.Range("B15:K16").FormatConditions
for ($i = 1; $i -le $fcs2.Count; $i++) {
    $fcs2.Item($i).ModifyAppliesToRange($ws2.Range("B15:L16"))
}

$wb.Application.Calculate()
